# Applies the crypto price/volume update described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '44.042.95'
$ws.Range("E2").Value = '  +0.98%  '
$ws.Range("D3").Value = '2.353.97'
$ws.Range("E3").Value = '  +0.10%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.695'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +5.88%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '241.69'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +3.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '76.30'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +4.24%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.632'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +21.12%  '
$ws.Range("E10").Value = '  +4.02%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '57.38'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.95%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '33.38'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +22.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '7.51'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +13.67%  '
$ws.Range("E14").Value = '  +1.72%  '
$ws.Range("D15").Value = '2.706.37'
$ws.Range("E15").Value = '  -0.03%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '16.82'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.929'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +6.19%  '
$ws.Range("D18").Value = '2.357.68'
$ws.Range("E18").Value = '  -0.69%  '
$ws.Range("D19").Value = '43.932.49'
$ws.Range("E19").Value = '  +0.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000103'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.34%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.68'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.27%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '77.73'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.76%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '262.69'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.81%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.00'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.54'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.43%  '
$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.81'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +18.11%  '
$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.63'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -4.85%  '
$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.94'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +7.51%  '
$ws.Range("E29").Value = '  +2.85%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '23.14'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +3.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '175.24'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.127'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -4.33%  '
$ws.Range("E33").Value = '  +5.28%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.40'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.58%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0762'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +9.23%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.45'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.50%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.82'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.39%  '
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("E39").Value = '  -2.37%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0283'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +7.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.221'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +22.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '19.28'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.62%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '9.13'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.91%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.107'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +11.20%  '
$ws.Range("E45").Value = '  -0.09%  '
$ws.Range("B46").Value = 'FTXToken'
$ws.Range("C46").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.69'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.64%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.53'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +11.65%  '
$ws.Range("B48").Value = 'TrustWalletToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.26'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.77%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '101.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.06%  '
$ws.Range("E50").Value = '  +2.15%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '56.55'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +11.31%  '
